# Apply updated "想去人数" (F column) values across the relevant sheets,
# matching the data refresh captured in the commit "Update gh-pages to
# output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 212
$ws1.Range("F5").Value  = 187
$ws1.Range("F9").Value  = 9942
$ws1.Range("F11").Value = 2868
$ws1.Range("F12").Value = 207
$ws1.Range("F13").Value = 2414
$ws1.Range("F14").Value = 2722
$ws1.Range("F17").Value = 2117
$ws1.Range("F18").Value = 44
$ws1.Range("F20").Value = 374
$ws1.Range("F22").Value = 103
$ws1.Range("F25").Value = 183
$ws1.Range("F30").Value = 127
$ws1.Range("F32").Value = 1951
$ws1.Range("F33").Value = 2887
$ws1.Range("F34").Value = 5
$ws1.Range("F36").Value = 1014
$ws1.Range("F39").Value = 1283
$ws1.Range("F41").Value = 94
$ws1.Range("F43").Value = 27

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 965
$ws3.Range("F4").Value = 118
$ws3.Range("F5").Value = 1924

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 965
$ws4.Range("F4").Value  = 118
$ws4.Range("F9").Value  = 187
$ws4.Range("F13").Value = 9942
$ws4.Range("F16").Value = 2869
$ws4.Range("F17").Value = 207
$ws4.Range("F18").Value = 2414
$ws4.Range("F19").Value = 2722
$ws4.Range("F21").Value = 2117
$ws4.Range("F22").Value = 44
$ws4.Range("F24").Value = 374
$ws4.Range("F28").Value = 183
$ws4.Range("F32").Value = 127
$ws4.Range("F34").Value = 1951
$ws4.Range("F36").Value = 2887
$ws4.Range("F37").Value = 1014
$ws4.Range("F44").Value = 1283
$ws4.Range("F47").Value = 27
